$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A49").Value = 46008
$ws.Range("B49").Value = 661
$ws.Range("C49").Value = 12
$ws.Range("D49").Value = 649

$ws.Range("A49:D49").Select()
